$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column B (shifts old B->C, C->D, D->E, E->F)
# This preserves the original best-fit column widths/styles for the shifted columns,
# matching how the author likely performed the edit in the Excel UI.
$ws.Columns.Item(2).Insert() | Out-Null

# Clear the stale "Best/Brother/Veg" block that used to live at row 15
# (its column positions changed completely, so clear first then rewrite).
$ws.Range("A15:E16").ClearContents() | Out-Null

# Write every cell to match the final desired content exactly.
$ws.Range("C12").Value = "Harry Cheung Wife"
$ws.Range("D25").Value = "Ellen "
$ws.Range("D26").Value = "Shirley Tsang"
$ws.Range("D22").Value = "Ah Chicken"
$ws.Range("D23").Value = "Bean Man"
$ws.Range("D24").Value = "Foo Kwai"
$ws.Range("B1").Value = "R1"
$ws.Range("B2").Value = "R2"
$ws.Range("B3").Value = "R3"
$ws.Range("B13").Value = "?Wong Ming"
$ws.Range("B14").Value = "?Wong Kei"
$ws.Range("A22").Value = "?Reconfirm"
$ws.Range("A1").Value = "**Lam Kei"
$ws.Range("C1").Value = "*Anthony "
$ws.Range("D1").Value = "*Francis"
$ws.Range("E1").Value = "JZ"
$ws.Range("A2").Value = "**KaYan"
$ws.Range("C2").Value = "*Lam Kuen"
$ws.Range("D2").Value = "*Small Tony"
$ws.Range("E2").Value = "Huayi"
$ws.Range("A3").Value = "Me"
$ws.Range("C3").Value = "Faye"
$ws.Range("D3").Value = "*Jonathan"
$ws.Range("E3").Value = "Denis"
$ws.Range("A4").Value = "In"
$ws.Range("B4").Value = "*Yellow"
$ws.Range("C4").Value = "Rohda"
$ws.Range("D4").Value = "JonathanW"
$ws.Range("E4").Value = "Gabriel"
$ws.Range("A5").Value = "Dad"
$ws.Range("B5").Value = "Ms Tong"
$ws.Range("C5").Value = "Patrick"
$ws.Range("D5").Value = "#Somingtat"
$ws.Range("E5").Value = "Brilly"
$ws.Range("A6").Value = "Mum"
$ws.Range("B6").Value = "Denis"
$ws.Range("C6").Value = "Leo"
$ws.Range("D6").Value = "SomingtatW"
$ws.Range("E6").Value = "Tim"
$ws.Range("A7").Value = "Sis"
$ws.Range("B7").Value = "Begger"
$ws.Range("C7").Value = "Fruit"
$ws.Range("D7").Value = "Jasper"
$ws.Range("E7").Value = "Natalie"
$ws.Range("B8").Value = "Sol Bread"
$ws.Range("C8").Value = "Chan Dan"
$ws.Range("D8").Value = "Tin Shing"
$ws.Range("E8").Value = "Hao"
$ws.Range("B9").Value = "Fai Chi"
$ws.Range("C9").Value = "Ho Kim Chin"
$ws.Range("D9").Value = "See Fu"
$ws.Range("E9").Value = "Sean"
$ws.Range("B10").Value = "Shum Kit"
$ws.Range("C10").Value = "Ho Ming"
$ws.Range("D10").Value = "Yvoone"
$ws.Range("E10").Value = "Arbinnav"
$ws.Range("B11").Value = "Zuey Tsui"
$ws.Range("C11").Value = "Harry Cheung"
$ws.Range("D11").Value = "Pui"
$ws.Range("E11").Value = "Denvendra"
$ws.Range("B12").Value = "Mak Wing"
$ws.Range("D12").Value = "Ocean"
$ws.Range("E12").Value = "Casey"
$ws.Range("A19").Value = "**Best"
$ws.Range("D19").Value = "Guanglei"
$ws.Range("E19").Value = "Billy"
$ws.Range("A20").Value = "*Brother"
$ws.Range("D20").Value = "Tim"
$ws.Range("E20").Value = "Arun"
$ws.Range("A21").Value = "#Veg"
$ws.Range("D21").Value = "Patrick"

# Clear any leftover cells in row 15/16 beyond E that may remain from the shift (safety net).
$ws.Range("F1:F30").ClearContents() | Out-Null

# Re-apply the explicit "no fill" style flag to the six cells that carry it in the target file.
$ws.Range("C1").Interior.Pattern = -4142
$ws.Range("D1").Interior.Pattern = -4142
$ws.Range("C2").Interior.Pattern = -4142
$ws.Range("D2").Interior.Pattern = -4142
$ws.Range("D3").Interior.Pattern = -4142
$ws.Range("B4").Interior.Pattern = -4142

# Adjust column widths: keep column B unformatted (no bestFit), and size the trailing columns.
$ws.Columns.Item(2).ColumnWidth = 9.0
$ws.Columns.Item(5).ColumnWidth = 11.166666666666666
$ws.Columns.Item(6).ColumnWidth = 12.666666666666666

# Restore the active selection cell shown in the saved workbook.
$ws.Range("D13").Select() | Out-Null

